$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row 11: n-Octano ---------------------------------------------
$ws.Range("A11").Value = "n-Octano"

# Formula column with subscripted "8" and "18" (C8H18)
$ws.Range("B11").Value = "C8H18"
$b11 = $ws.Range("B11")
$b11.Characters(2,1).Font.Subscript = $true
$b11.Characters(3,1).Font.Name = "Calibri"
$b11.Characters(4,2).Font.Subscript = $true

$ws.Range("C11").Value = 4.05075
$ws.Range("D11").Value = 1356.36
$ws.Range("E11").Value = 209.635
$ws.Range("F11").Value = 295.4
$ws.Range("G11").Value = 24.82463

# --- Formatting ---------------------------------------------------------
# A11:E11 -> wrap text, vertically centered, no border (matches style used
# elsewhere in the sheet minus the box border)
$rngAE = $ws.Range("A11:E11")
$rngAE.VerticalAlignment = -4108   # xlCenter
$rngAE.WrapText = $true

# F11 -> same alignment, plus thin left/right border
$rngF = $ws.Range("F11")
$rngF.VerticalAlignment = -4108    # xlCenter
$rngF.WrapText = $true
$rngF.Borders.Item(7).LineStyle = 1   # xlEdgeLeft
$rngF.Borders.Item(10).LineStyle = 1  # xlEdgeRight

# --- Workbook / window bookkeeping -------------------------------------
$win = $excel.ActiveWindow
$win.Left = 4128
$win.Top = 1104
$win.Width = 17280
$win.Height = 8964

# Match the final selection recorded in the workbook
$ws.Range("G11").Select()
